# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Sat Jul 13 23:59:34 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "59.209.71"
$ws.Range("E2").Value = "  +2.25%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "3.175.75"
$ws.Range("E3").Value = "  +1.34%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.00%  "

# --- Row 5: BNB ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.35"

# --- Row 6: Solana ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.94"
$ws.Range("E6").Value = "  +1.75%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.02%  "

# --- Row 8: XRP ---
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  +10.92%  "

# --- Row 9: Toncoin ---
$ws.Range("E9").Value = "  -0.86%  "

# --- Row 10: Cardano ---
$ws.Range("E10").Value = "  +5.93%  "

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = "  +4.34%  "

# --- Row 12: becomes TRON (was WrappedliquidstakedEther2.0) ---
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.139"
$ws.Range("E12").Value = "  +1.59%  "

# --- Row 13: becomes WrappedliquidstakedEther2.0 (was TRON) ---
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.724.65"
$ws.Range("E13").Value = "  +1.71%  "

# --- Row 14: Avalanche ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.87"
$ws.Range("E14").Value = "  -0.37%  "

# --- Row 15: ShibaInu ---
$ws.Range("E15").Value = "  +3.77%  "

# --- Row 16: WrappedBTC ---
$ws.Range("D16").Value = "59.229.99"
$ws.Range("E16").Value = "  +2.14%  "

# --- Row 17: becomes WrappedEther (was Polkadot) ---
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.167.77"
$ws.Range("E17").Value = "  +1.17%  "

# --- Row 18: becomes Polkadot (was WrappedEther) ---
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.24"
$ws.Range("E18").Value = "  +2.35%  "

# --- Row 19: Chainlink ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("E19").Value = "  +2.24%  "

# --- Row 20: Uniswap ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.16"
$ws.Range("E20").Value = "  +0.60%  "

# --- Row 21: BitcoinCash ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.11"
$ws.Range("E21").Value = "  +1.99%  "

# --- Row 22: Dai ---
$ws.Range("E22").Value = "  +0.06%  "

# --- Row 23: Polygon ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.531"
$ws.Range("E23").Value = "  +4.77%  "

# --- Row 24: Litecoin ---
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.69"
$ws.Range("E24").Value = "  +0.54%  "

# --- Row 25: Kaspa ---
$ws.Range("E25").Value = "  -0.62%  "

# --- Row 26: Binance-PegBSC-USD ---
$ws.Range("E26").Value = "  +0.06%  "

# --- Row 27: InternetComputer(DFINITY) ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.44"
$ws.Range("E27").Value = "  +15.28%  "

# --- Row 28: PEPE ---
$ws.Range("D28").Value = "0.0₃0871"
$ws.Range("E28").Value = "  +0.27%  "

# --- Row 29: EthereumClassic ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.41"
$ws.Range("E29").Value = "  +4.39%  "

# --- Row 30: PancakeSwap ---
$ws.Range("E30").Value = "  +0.45%  "

# --- Row 31: RenderToken ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.02"
$ws.Range("E31").Value = "  -1.50%  "

# --- Row 32: NEARProtocol ---
$ws.Range("E32").Value = "  +0.49%  "

# --- Row 33: Fetch.AI ---
$ws.Range("E33").Value = "  -2.06%  "

# --- Row 34: Aptos ---
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.34"
$ws.Range("E34").Value = "  +3.95%  "

# --- Row 35: Monero ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.26"
$ws.Range("E35").Value = "  -1.38%  "

# --- Row 36: ImmutableX ---
$ws.Range("E36").Value = "  +3.40%  "

# --- Row 37: Hedera ---
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0713"
$ws.Range("E37").Value = "  +5.95%  "

# --- Row 38: EnergySwap ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.37"
$ws.Range("E38").Value = "  -0.87%  "

# --- Row 39: Maker ---
$ws.Range("D39").Value = "2.709.28"
$ws.Range("E39").Value = "  +7.67%  "

# --- Row 40: Stacks ---
$ws.Range("E40").Value = "  +1.36%  "

# --- Row 41: Filecoin ---
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.28"
$ws.Range("E41").Value = "  +3.90%  "

# --- Row 42: VeChain ---
$ws.Range("E42").Value = "  +8.22%  "

# --- Row 43: Mantle ---
$ws.Range("E43").Value = "  +3.38%  "

# --- Row 44: OKB ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.12"
$ws.Range("E44").Value = "  +3.45%  "

# --- Row 45: FirstDigitalUSD ---
$ws.Range("E45").Value = "  -0.02%  "

# --- Row 46: RenzoRestakedETH ---
$ws.Range("D46").Value = "3.216.75"
$ws.Range("E46").Value = "  +1.27%  "

# --- Row 47: Stellar ---
$ws.Range("E47").Value = "  +12.23%  "

# --- Row 48: ONDO ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.984"
$ws.Range("E48").Value = "  +0.04%  "

# --- Row 49: Cosmos ---
$ws.Range("E49").Value = "  +0.89%  "

# --- Row 50: InjectiveProtocol ---
$ws.Range("E50").Value = "  +2.47%  "

# --- Row 51: SuiNetwork ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.758"
$ws.Range("E51").Value = "  +1.34%  "
